# Formed the consolidated report
# The "Absent" column (H) should reflect whether there was a "Real" (column E)
# attendance recorded for the day: if Real == 0 then Absent = 1, otherwise Absent = 0.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

# Data rows start at row 3 (rows 1-2 are headers / student info)
for ($r = 3; $r -le $lastRow; $r++) {
    $realCell = $ws.Cells.Item($r, 5)   # Column E = Real
    $absentCell = $ws.Cells.Item($r, 8) # Column H = Absent

    $realValue = $realCell.Value()
    if ($null -eq $realValue -or $realValue -eq "" -or $realValue -eq 0) {
        $absentCell.Value = 1
    } else {
        $absentCell.Value = 0
    }
}
